$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1709
$ws.Range("F4").Value = 200
$ws.Range("F5").Value = 1079
$ws.Range("F6").Value = 29
$ws.Range("F7").Value = 132
$ws.Range("F8").Value = 1387
$ws.Range("F11").Value = 397
$ws.Range("F13").Value = 69
$ws.Range("F15").Value = 430
$ws.Range("F16").Value = 457
$ws.Range("F17").Value = 123
$ws.Range("F18").Value = 22
$ws.Range("F19").Value = 441
$ws.Range("F22").Value = 36
$ws.Range("F26").Value = 166
$ws.Range("F28").Value = 91
$ws.Range("F29").Value = 351
$ws.Range("F31").Value = 43
$ws.Range("F32").Value = 21
$ws.Range("F36").Value = 230

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 56
$ws.Range("E5").Value = "2024.02.04 19:30-02.04 21:30"
$ws.Range("F5").Value = 587
$ws.Range("E6").Value = "2024.02.04 19:30-02.04 21:30"
$ws.Range("F6").Value = 587
$ws.Range("F12").Value = 265
$ws.Range("F15").Value = 285
$ws.Range("F16").Value = 285
$ws.Range("F18").Value = 7
$ws.Range("F19").Value = 919
$ws.Range("F22").Value = 592
$ws.Range("F24").Value = 20
$ws.Range("F26").Value = 207
$ws.Range("F28").Value = 12
$ws.Range("F33").Value = 6

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1749
$ws.Range("F5").Value = 2041
$ws.Range("F6").Value = 2247
$ws.Range("F7").Value = 885
$ws.Range("F10").Value = 1080
$ws.Range("F11").Value = 225

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1749
$ws.Range("F5").Value = 2041
$ws.Range("F6").Value = 2247
$ws.Range("F8").Value = 1709
$ws.Range("F10").Value = 56
$ws.Range("F11").Value = 885
$ws.Range("F12").Value = 1080
$ws.Range("F13").Value = 200
$ws.Range("F14").Value = 225
$ws.Range("F17").Value = 1079
$ws.Range("F18").Value = 132
$ws.Range("E19").Value = "2024.02.04 19:30-02.04 21:30"
$ws.Range("F19").Value = 587
$ws.Range("F22").Value = 397
$ws.Range("F24").Value = 69
$ws.Range("F26").Value = 430
$ws.Range("F27").Value = 457
$ws.Range("F28").Value = 123
$ws.Range("F29").Value = 22
$ws.Range("F30").Value = 442
$ws.Range("F35").Value = 166
$ws.Range("F36").Value = 91
$ws.Range("F38").Value = 351
$ws.Range("F40").Value = 285
$ws.Range("F41").Value = 43
$ws.Range("F43").Value = 20
$ws.Range("F44").Value = 207
$ws.Range("F49").Value = 230
